$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer_churn_Filter")

# Update the credit card segmentation label to the shorter form.
$ws.Range("C2").Value = "Credit Cards"

# Reflect the cell the user left selected after making the edit.
$ws.Activate()
$ws.Range("C2").Select()
